$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text (non-numeric-looking) cells: safe to assign directly ---
$ws.Range("B4").Value = "2024.12.01. 월간"
$ws.Range("B7").Value = "2025년 01월 10일 01시 33분 24초"

# --- Numeric-looking values that must stay stored as TEXT (as in the source
#     workbook, which uses inlineStr / text cells even for numbers).
#     Assigning ".Value" directly on these cells would make Excel parse the
#     string and store it as a real number, changing cell type + losing the
#     original styling. Instead: write the text into a helper cell that is
#     explicitly formatted as Text ("@"), copy it, and paste-special VALUES
#     ONLY into the target cell so the target keeps its original style/format
#     and only the stored value/type changes. The helper cell is cleared
#     (format included) afterwards so it leaves no trace on the sheet.

$helper = $ws.Range("Z1")

function Set-TextValue($cellAddress, $text) {
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($cellAddress).PasteSpecial(-4163)
}

$pairs = @(
    @("C10","78"),    @("D10","29"),
    @("C11","191"),   @("D11","71"),
    @("D14","0.37"),
    @("C15","5"),     @("D15","1.86"),
    @("C16","29"),    @("D16","10.78"),
    @("C17","117"),   @("D17","43.49"),
    @("C18","16"),    @("D18","5.95"),
    @("C19","21"),    @("D19","7.81"),
    @("C20","9"),     @("D20","3.35"),
    @("C21","9"),     @("D21","3.35"),
    @("C22","7"),     @("D22","2.6"),
    @("C23","7"),     @("D23","2.6"),
    @("C24","5"),     @("D24","1.86"),
    @("C25","8"),     @("D25","2.97"),
    @("C26","2"),     @("D26","0.74"),
    @("D27","4.09"),
    @("C28","3"),     @("D28","1.12"),
    @("C29","4"),     @("D29","1.49"),
    @("C30","3"),     @("D30","1.12"),
    @("D31","1.49"),
    @("C32","3"),     @("D32","1.12"),
    @("C33","5"),     @("D33","1.86")
)

foreach ($pair in $pairs) {
    Set-TextValue $pair[0] $pair[1]
}

$helper.Clear()

Write-Host "Applied $($pairs.Count) numeric-text updates plus 2 header text updates."
